$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 28.09534366666667
$ws.Range("H2").Value = 84.28603100000001
$ws.Range("I2").Value = 0.6431807885786103
$ws.Range("J2").Value = 0.6438169333507339
$ws.Range("M2").Value = 0.013841
$ws.Range("N2").Value = 0.041523
$ws.Range("O2").Value = 0.001379248562145083
$ws.Range("P2").Value = 0.001433331793786624
$ws.Range("Q2").Value = 0.3888676516903333
$ws.Range("R2").Value = 3.499808865213
$ws.Range("S2").Value = 0.0008871061778463887
$ws.Range("T2").Value = 0.0009228032799498108

# Row 3
$ws.Range("G3").Value = 28.09534366666667
$ws.Range("H3").Value = 84.28603100000001
$ws.Range("I3").Value = 0.6431807885786103
$ws.Range("J3").Value = 0.6438169333507339
$ws.Range("M3").Value = 3.248047
$ws.Range("N3").Value = 9.744140999999999
$ws.Range("O3").Value = 0.3236662202535691
$ws.Range("P3").Value = 0.3363578522370683
$ws.Range("Q3").Value = 91.25499671048566
$ws.Range("R3").Value = 821.294970394371
$ws.Range("S3").Value = 0.2081758947789487
$ws.Range("T3").Value = 0.2165528809357086

# Row 4
$ws.Range("G4").Value = 28.09534366666667
$ws.Range("H4").Value = 84.28603100000001
$ws.Range("I4").Value = 0.6431807885786103
$ws.Range("J4").Value = 0.6438169333507339
$ws.Range("M4").Value = 1.1359575
$ws.Range("N4").Value = 2.271915
$ws.Range("O4").Value = 0.1131975831611099
$ws.Range("P4").Value = 0.07842419869182714
$ws.Range("Q4").Value = 31.9151163532275
$ws.Range("R4").Value = 191.490698119365
$ws.Range("S4").Value = 0.07280651080275551
$ws.Range("T4").Value = 0.05049082710226079

# Row 5
$ws.Range("G5").Value = 28.09534366666667
$ws.Range("H5").Value = 84.28603100000001
$ws.Range("I5").Value = 0.6431807885786103
$ws.Range("J5").Value = 0.6438169333507339
$ws.Range("M5").Value = 5.637329
$ws.Range("N5").Value = 16.911987
$ws.Range("O5").Value = 0.5617569480231759
$ws.Range("P5").Value = 0.5837846172773179
$ws.Range("Q5").Value = 158.3826956170664
$ws.Range("R5").Value = 1425.444260553597
$ws.Range("S5").Value = 0.3613112768190597
$ws.Range("T5").Value = 0.3758504220328147

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.16980166666667
$ws.Range("H6").Value = 45.509405
$ws.Range("I6").Value = 0.347279076358968
$ws.Range("J6").Value = 0.3476225564081497
$ws.Range("M6").Value = 0.013841
$ws.Range("N6").Value = 0.041523
$ws.Range("O6").Value = 0.001379248562145083
$ws.Range("P6").Value = 0.001433331793786624
$ws.Range("Q6").Value = 0.2099652248683333
$ws.Range("R6").Value = 1.889687023815
$ws.Range("S6").Value = 0.0004789841667311791
$ws.Range("T6").Value = 0.000498258462337185

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.16980166666667
$ws.Range("H7").Value = 45.509405
$ws.Range("I7").Value = 0.347279076358968
$ws.Range("J7").Value = 0.3476225564081497
$ws.Range("M7").Value = 3.248047
$ws.Range("N7").Value = 9.744140999999999
$ws.Range("O7").Value = 0.3236662202535691
$ws.Range("P7").Value = 0.3363578522370683
$ws.Range("Q7").Value = 49.27222879401166
$ws.Range("R7").Value = 443.450059146105
$ws.Range("S7").Value = 0.1124025060182578
$ws.Range("T7").Value = 0.1169255764626043

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 15.16980166666667
$ws.Range("H8").Value = 45.509405
$ws.Range("I8").Value = 0.347279076358968
$ws.Range("J8").Value = 0.3476225564081497
$ws.Range("M8").Value = 1.1359575
$ws.Range("N8").Value = 2.271915
$ws.Range("O8").Value = 0.1131975831611099
$ws.Range("P8").Value = 0.07842419869182714
$ws.Range("Q8").Value = 17.2322499767625
$ws.Range("R8").Value = 103.393499860575
$ws.Range("S8").Value = 0.03931115212625774
$ws.Range("T8").Value = 0.02726202043351362

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 15.16980166666667
$ws.Range("H9").Value = 45.509405
$ws.Range("I9").Value = 0.347279076358968
$ws.Range("J9").Value = 0.3476225564081497
$ws.Range("M9").Value = 5.637329
$ws.Range("N9").Value = 16.911987
$ws.Range("O9").Value = 0.5617569480231759
$ws.Range("P9").Value = 0.5837846172773179
$ws.Range("Q9").Value = 85.51716285974834
$ws.Range("R9").Value = 769.654465737735
$ws.Range("S9").Value = 0.1950864340477213
$ws.Range("T9").Value = 0.2029367010496945

# Row 10
$ws.Range("G10").Value = 0.129484
$ws.Range("H10").Value = 0.258968
$ws.Range("I10").Value = 0.002964249956021043
$ws.Range("J10").Value = 0.001978121185893458
$ws.Range("M10").Value = 0.013841
$ws.Range("N10").Value = 0.041523
$ws.Range("O10").Value = 0.001379248562145083
$ws.Range("P10").Value = 0.001433331793786624
$ws.Range("Q10").Value = 0.001792188044
$ws.Range("R10").Value = 0.010753128264
$ws.Range("S10").Value = 0.000004088437489680648
$ws.Range("T10").Value = 0.000002835303987703995

# Row 11
$ws.Range("G11").Value = 0.129484
$ws.Range("H11").Value = 0.258968
$ws.Range("I11").Value = 0.002964249956021043
$ws.Range("J11").Value = 0.001978121185893458
$ws.Range("M11").Value = 3.248047
$ws.Range("N11").Value = 9.744140999999999
$ws.Range("O11").Value = 0.3236662202535691
$ws.Range("P11").Value = 0.3363578522370683
$ws.Range("Q11").Value = 0.4205701177479999
$ws.Range("R11").Value = 2.523420706487999
$ws.Range("S11").Value = 0.0009594275791521393
$ws.Range("T11").Value = 0.0006653565935517662

# Row 12
$ws.Range("G12").Value = 0.129484
$ws.Range("H12").Value = 0.258968
$ws.Range("I12").Value = 0.002964249956021043
$ws.Range("J12").Value = 0.001978121185893458
$ws.Range("M12").Value = 1.1359575
$ws.Range("N12").Value = 2.271915
$ws.Range("O12").Value = 0.1131975831611099
$ws.Range("P12").Value = 0.07842419869182714
$ws.Range("Q12").Value = 0.14708832093
$ws.Range("R12").Value = 0.58835328372
$ws.Range("S12").Value = 0.0003355459309070085
$ws.Range("T12").Value = 0.0001551325689190213

# Row 13
$ws.Range("G13").Value = 0.129484
$ws.Range("H13").Value = 0.258968
$ws.Range("I13").Value = 0.002964249956021043
$ws.Range("J13").Value = 0.001978121185893458
$ws.Range("M13").Value = 5.637329
$ws.Range("N13").Value = 16.911987
$ws.Range("O13").Value = 0.5617569480231759
$ws.Range("P13").Value = 0.5837846172773179
$ws.Range("Q13").Value = 0.729943908236
$ws.Range("R13").Value = 4.379663449415999
$ws.Range("S13").Value = 0.001665188008472214
$ws.Range("T13").Value = 0.001154796719434967

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.287247
$ws.Range("H14").Value = 0.861741
$ws.Range("I14").Value = 0.006575885106400611
$ws.Range("J14").Value = 0.006582389055223097
$ws.Range("M14").Value = 0.013841
$ws.Range("N14").Value = 0.041523
$ws.Range("O14").Value = 0.001379248562145083
$ws.Range("P14").Value = 0.001433331793786624
$ws.Range("Q14").Value = 0.003975785726999999
$ws.Range("R14").Value = 0.035782071543
$ws.Range("S14").Value = 0.000009069780077834306
$ws.Range("T14").Value = 0.000009434747511924363

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.287247
$ws.Range("H15").Value = 0.861741
$ws.Range("I15").Value = 0.006575885106400611
$ws.Range("J15").Value = 0.006582389055223097
$ws.Range("M15").Value = 3.248047
$ws.Range("N15").Value = 9.744140999999999
$ws.Range("O15").Value = 0.3236662202535691
$ws.Range("P15").Value = 0.3363578522370683
$ws.Range("Q15").Value = 0.9329917566089998
$ws.Range("R15").Value = 8.396925809480999
$ws.Range("S15").Value = 0.002128391877210424
$ws.Range("T15").Value = 0.002214038245203626

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.287247
$ws.Range("H16").Value = 0.861741
$ws.Range("I16").Value = 0.006575885106400611
$ws.Range("J16").Value = 0.006582389055223097
$ws.Range("M16").Value = 1.1359575
$ws.Range("N16").Value = 2.271915
$ws.Range("O16").Value = 0.1131975831611099
$ws.Range("P16").Value = 0.07842419869182714
$ws.Range("Q16").Value = 0.3263003840024999
$ws.Range("R16").Value = 1.957802304015
$ws.Range("S16").Value = 0.0007443743011896875
$ws.Range("T16").Value = 0.0005162185871337245

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.287247
$ws.Range("H17").Value = 0.861741
$ws.Range("I17").Value = 0.006575885106400611
$ws.Range("J17").Value = 0.006582389055223097
$ws.Range("M17").Value = 5.637329
$ws.Range("N17").Value = 16.911987
$ws.Range("O17").Value = 0.5617569480231759
$ws.Range("P17").Value = 0.5837846172773179
$ws.Range("Q17").Value = 1.619305843263
$ws.Range("R17").Value = 14.573752589367
$ws.Range("S17").Value = 0.003694049147922664
$ws.Range("T17").Value = 0.003842697475373822

